# Update the public EPEX Spot prices workbook.
#
# 1) "Prix Spot" sheet: a new day ("22-dec") is inserted as a new column
#    right before the "01-oct." column (currently column ES). All the
#    "01-oct." ... "31-oct." columns shift one column to the right
#    (ES:FW -> ET:FX). The brand-new column has no data yet, so every
#    data row (2-25) gets the usual placeholder "-".
# 2) "Gaz" sheet: two new daily rows are appended (2025-12-20, 2025-12-21),
#    carrying forward the last known price (26.9).
# 3) "CO2" sheet: two new daily rows are appended (2025-12-20, 2025-12-21),
#    carrying forward the last known price (84.54000000000001).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Prix Spot : insert a new "22-dec" column before the "01-oct." column
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Inserting the entire column shifts ES:FW (01-oct. .. 31-oct.) to ET:FX
# and copies the formatting of the former ES column into the new one.
$wsSpot.Range("ES1").EntireColumn.Insert()

# Header for the newly inserted column.
$wsSpot.Range("ES1").Value = "22-dec"

# No data is available yet for this new day, so fill it with "-" like
# every other not-yet-known day in the sheet.
$wsSpot.Range("ES2:ES25").Value = "-"

# ---------------------------------------------------------------------
# 2) Gaz : append 2025-12-20 and 2025-12-21 (same price as 2025-12-19)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force text (the date column stores plain text like "2025-06-16", not
# real Excel dates) before assigning the value, otherwise Excel would
# auto-convert the ISO-like string into a date serial number. Resetting
# the Style back to "Normal" afterwards drops the temporary Text number
# format again so the cell ends up with the same (default) style as all
# the other date cells in the column.
$wsGaz.Range("A177").NumberFormat = "@"
$wsGaz.Range("A177").Value = "2025-12-20"
$wsGaz.Range("A177").Style = "Normal"
$wsGaz.Range("B177").Value = 26.9

$wsGaz.Range("A178").NumberFormat = "@"
$wsGaz.Range("A178").Value = "2025-12-21"
$wsGaz.Range("A178").Style = "Normal"
$wsGaz.Range("B178").Value = 26.9

# ---------------------------------------------------------------------
# 3) CO2 : append 2025-12-20 and 2025-12-21 (same price as 2025-12-19)
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A178").NumberFormat = "@"
$wsCO2.Range("A178").Value = "2025-12-20"
$wsCO2.Range("A178").Style = "Normal"
$wsCO2.Range("B178").Value = 84.54000000000001

$wsCO2.Range("A179").NumberFormat = "@"
$wsCO2.Range("A179").Value = "2025-12-21"
$wsCO2.Range("A179").Style = "Normal"
$wsCO2.Range("B179").Value = 84.54000000000001
